# Updated cryptos list with refreshed Price (D) and Volume(1h) (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.717.37"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +1.20%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.873.61"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +1.45%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.14%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "332.04"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +3.44%  "

$ws.Range("E6").Value = "  +0.14%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4743"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +6.52%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3954"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +3.33%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "47.56"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -3.69%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.08034"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +2.81%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.021"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +0.95%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "21.82"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +2.03%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.904.23"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +3.71%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.954"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +2.18%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.140"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +0.96%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +0.20%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "87.14"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +2.51%  "

$ws.Range("E18").Value = "  +1.89%  "

$ws.Range("E19").Value = "  +2.50%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.27"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +1.95%  "

$ws.Range("E21").Value = "  +0.27%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "27.737.88"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +1.33%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.490"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("E24").Value = "  +2.50%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.300"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +1.72%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.105.66"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +2.35%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "156.42"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +3.43%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "20.20"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +4.92%  "

$ws.Range("E29").Value = "  +3.31%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.559"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +2.10%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "122.36"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +2.22%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.9685"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +5.08%  "

$ws.Range("E33").Value = "  +3.13%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.451"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -1.33%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.632"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +1.13%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.293"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +1.60%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.06107"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +3.15%  "

$ws.Range("E38").Value = "  +2.41%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.223"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.41%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "8.158"
$cell.ClearFormats()

$ws.Range("E41").Value = "  +0.13%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.5988"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +1.78%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1908"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +3.65%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "10.25"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +0.62%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.254"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +0.21%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5687"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +0.86%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "12.29"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +1.26%  "

$ws.Range("E48").Value = "  +1.65%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.930"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.90%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.06817"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("E51").Value = "  +10.61%  "
